$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.29903
$ws.Range("H2").Value = 6.89709
$ws.Range("I2").Value = 0.04075801785348079
$ws.Range("J2").Value = 0.04075801785348079
$ws.Range("M2").Value = 29.80827733333334
$ws.Range("N2").Value = 89.42483200000001
$ws.Range("O2").Value = 0.4866975737940222
$ws.Range("P2").Value = 0.4866975737940221
$ws.Range("Q2").Value = 68.53012383765335
$ws.Range("R2").Value = 616.7711145388801
$ws.Range("S2").Value = 0.01983682840194254
$ws.Range("T2").Value = 0.01983682840194254

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.29903
$ws.Range("H3").Value = 6.89709
$ws.Range("I3").Value = 0.04075801785348079
$ws.Range("J3").Value = 0.04075801785348079
$ws.Range("O3").Value = 0.437868100938039
$ws.Range("P3").Value = 0.437868100938039
$ws.Range("Q3").Value = 61.65462249569666
$ws.Range("R3").Value = 554.8916024612701
$ws.Range("S3").Value = 0.01784663587550232
$ws.Range("T3").Value = 0.01784663587550232

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.29903
$ws.Range("H4").Value = 6.89709
$ws.Range("I4").Value = 0.04075801785348079
$ws.Range("J4").Value = 0.04075801785348079
$ws.Range("M4").Value = 4.620050333333332
$ws.Range("O4").Value = 0.07543432526793886
$ws.Range("P4").Value = 0.07543432526793886
$ws.Range("Q4").Value = 10.62163431784333
$ws.Range("R4").Value = 95.59470886058999
$ws.Range("S4").Value = 0.003074553576035929
$ws.Range("T4").Value = 0.003074553576035929

$ws.Range("I5").Value = 0.3949230674234065
$ws.Range("J5").Value = 0.3949230674234066
$ws.Range("M5").Value = 29.80827733333334
$ws.Range("N5").Value = 89.42483200000001
$ws.Range("O5").Value = 0.4866975737940222
$ws.Range("P5").Value = 0.4866975737940221
$ws.Range("Q5").Value = 664.0196982631396
$ws.Range("R5").Value = 5976.177284368257
$ws.Range("S5").Value = 0.192208098750265
$ws.Range("T5").Value = 0.192208098750265

$ws.Range("I6").Value = 0.3949230674234065
$ws.Range("J6").Value = 0.3949230674234066
$ws.Range("O6").Value = 0.437868100938039
$ws.Range("P6").Value = 0.437868100938039
$ws.Range("S6").Value = 0.1729242135493122
$ws.Range("T6").Value = 0.1729242135493122

$ws.Range("I7").Value = 0.3949230674234065
$ws.Range("J7").Value = 0.3949230674234066
$ws.Range("M7").Value = 4.620050333333332
$ws.Range("O7").Value = 0.07543432526793886
$ws.Range("P7").Value = 0.07543432526793886
$ws.Range("R7").Value = 926.260835067758
$ws.Range("S7").Value = 0.0297907551238294
$ws.Range("T7").Value = 0.0297907551238294

$ws.Range("I8").Value = 0.5643189147231126
$ws.Range("J8").Value = 0.5643189147231126
$ws.Range("M8").Value = 29.80827733333334
$ws.Range("N8").Value = 89.42483200000001
$ws.Range("O8").Value = 0.4866975737940222
$ws.Range("P8").Value = 0.4866975737940221
$ws.Range("Q8").Value = 948.8401827814188
$ws.Range("R8").Value = 8539.561645032769
$ws.Range("S8").Value = 0.2746526466418146
$ws.Range("T8").Value = 0.2746526466418145

$ws.Range("I9").Value = 0.5643189147231126
$ws.Range("J9").Value = 0.5643189147231126
$ws.Range("O9").Value = 0.437868100938039
$ws.Range("P9").Value = 0.437868100938039
$ws.Range("S9").Value = 0.2470972515132245
$ws.Range("T9").Value = 0.2470972515132245

$ws.Range("I10").Value = 0.5643189147231126
$ws.Range("J10").Value = 0.5643189147231126
$ws.Range("M10").Value = 4.620050333333332
$ws.Range("O10").Value = 0.07543432526793886
$ws.Range("P10").Value = 0.07543432526793886
$ws.Range("Q10").Value = 147.0628226421276
$ws.Range("S10").Value = 0.04256901656807353
$ws.Range("T10").Value = 0.04256901656807353
